# Append a new quarterly period (column BB) to the year-over-year series:
# - BB1 gets the new period's date header (formatted like the other date
#   headers in row 1, i.e. copied from BA1's style).
# - BB3:BB21 repeat the most recent known YoY figure (column BA) as a
#   placeholder for the new period, matching the existing data rows.
# Row 2 and row 22 only ever contained column A, so they are untouched
# apart from their "spans" bookkeeping attribute (handled automatically by
# the engine once the sheet's used range grows to column BB).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell BB1 -------------------------------------------------
# Copy formatting (number format / font / border) from BA1, then set the
# new date serial value.
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("BB1").Value = 45986

# --- Data rows 3-21: mirror column BA into the new column BB --------
$lastRow = 21
for ($r = 3; $r -le $lastRow; $r++) {
    $srcCell = $ws.Range("BA$r")
    $srcValue = $srcCell.Value2()
    if ($null -ne $srcValue) {
        $ws.Range("BB$r").Value = $srcValue
    }
}
